$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (Wins / Losses / Ties) in columns AD, AE, AF of row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Reuse the existing header formatting (bold, centered, bordered) from
# column A's header cell, instead of creating a brand-new style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Team record columns for every data row (2-60): Wins=100, Losses=62, Ties=0.
for ($row = 2; $row -le 60; $row++) {
    $ws.Cells.Item($row, 30).Value = 100
    $ws.Cells.Item($row, 31).Value = 62
    $ws.Cells.Item($row, 32).Value = 0
}
